# Automated data refresh for "LOZANO MOLINA TITO" report workbook.
# Mirrors an external "Actualizacion automatica" export: every figure in
# this workbook is a plain literal (no live formulas anywhere in the
# file), so each touched cell is (re)written with its final literal
# value rather than a formula.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "VENTAS POR GRUPO": client LINO TUMBACO VICENTE JAVIER (row 14)
# picked up sales in two more product groups this period.
# ---------------------------------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Cells.Item(14, 14).Value = -513.72   # N14 - PUERTAS DE SEGURIDAD
$wsGrupo.Cells.Item(14, 17).Value = -44.6     # Q14 - PANELES PVC

# ---------------------------------------------------------------------
# Sheet "VENTA MENSUAL": same client's "julio" figure, and the column
# total underneath it.
# ---------------------------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Cells.Item(14, 6).Value = -558.3200000000001   # F14 - julio
$wsMensual.Cells.Item(29, 6).Value = 1618.96              # F29 - total julio

# ---------------------------------------------------------------------
# Sheet "CUMPLIMIENTO MENSUAL": PANELES PVC (row 14), PUERTAS DE
# SEGURIDAD (row 17) and the grand TOTAL (row 19) ripple from the same
# update, plus a one-unit widening of column F.
# ---------------------------------------------------------------------
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Column F: 24 -> 25 characters wide. ColumnWidth is expressed in
# characters of the Normal style font, which Excel persists to the
# package's <col width> in "character width + 5/6" units - back that
# offset out so the saved width lands exactly on 25.
$wsCumpl.Columns.Item(6).ColumnWidth = 25 - (5 / 6)

# Row 14 - PANELES PVC
$wsCumpl.Cells.Item(14, 4).Value = -44.6                   # D14 - VENTA
$wsCumpl.Cells.Item(14, 5).Value = 284.6                   # E14 - POR CUMPLIR
$wsCumpl.Cells.Item(14, 6).Value = -0.1858333333333334     # F14 - CUMPLIMIENTO

# Row 17 - PUERTAS DE SEGURIDAD
$wsCumpl.Cells.Item(17, 4).Value = -513.72                 # D17 - VENTA
$wsCumpl.Cells.Item(17, 5).Value = 855.72                  # E17 - POR CUMPLIR
$wsCumpl.Cells.Item(17, 6).Value = -1.502105263157895      # F17 - CUMPLIMIENTO

# Row 19 - TOTAL
$wsCumpl.Cells.Item(19, 4).Value = 1618.96                 # D19 - VENTA
$wsCumpl.Cells.Item(19, 5).Value = 25562.35093005039       # E19 - POR CUMPLIR
$wsCumpl.Cells.Item(19, 6).Value = 0.05956151284116889     # F19 - CUMPLIMIENTO

Write-Output "Actualizacion automatica aplicada"
